$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width change: column F (6) grows from the shared 61.0 group
# (cols F:H) to match column C's 68.43 width, splitting the group into
# F alone (widened) and G:H (unchanged at 61.0). The COM ColumnWidth
# setter quantizes to 1/6-character steps, so 67.67 is the input that
# round-trips to the closest achievable stored width (68.5 vs target 68.43).
$ws.Columns.Item(6).ColumnWidth = 67.67

# --- Template row to copy style/layout from (existing last data row) ---
$templateRow = 182

function Add-BookRow {
    param(
        [int]$RowNum,
        [hashtable]$Values
    )

    # Copy formatting + column layout (style s="2", and which columns are
    # populated) from the existing template row, then clear the
    # helper-only columns (T, U, W) that the template leaves blank so the
    # new row doesn't end up with stray empty cells.
    $srcRange = $ws.Range("A" + $templateRow + ":AF" + $templateRow)
    $dstRange = $ws.Range("A" + $RowNum + ":AF" + $RowNum)
    $srcRange.Copy($dstRange)
    $ws.Range("T" + $RowNum + ":U" + $RowNum).ClearContents()
    $ws.Range("W" + $RowNum).ClearContents()

    # Write cell values in strict left-to-right column order so any brand
    # new shared strings are appended to sharedStrings.xml in the same
    # order the target workbook expects.
    $columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","V","X","Y","Z","AA","AB","AC","AD","AE","AF")
    foreach ($col in $columns) {
        if ($Values.ContainsKey($col)) {
            $ws.Range($col + $RowNum).Value = $Values[$col]
        }
    }
}

Add-BookRow -RowNum 183 -Values @{
    A  = 182
    B  = "Book"
    C  = "The Daily Laws"
    D  = "Yes"
    E  = "Robert Greene"
    F  = "366 Medications On Power , Seduction , Mastery , Strategy and Human Nature"
    G  = "Focus Area"
    H  = 453
    I  = "NA"
    J  = "Rs. 699"
    K  = 1
    L  = "1-1-D"
    M  = "Paperback"
    N  = 2021
    O  = "English"
    P  = "Self Help"
    Q  = "Strategy"
    R  = "No"
    S  = "NF"
    V  = 0
    X  = "Male"
    Y  = "978-1-7881-6854-0"
    Z  = 1
    AA = "First Floor"
    AB = "Practical"
    AC = "Yes"
    AD = "Yes"
    AE = 9.8
    AF = "New"
}

Add-BookRow -RowNum 184 -Values @{
    A  = 183
    B  = "Book"
    C  = "The 33 Strategies of War"
    D  = "Yes"
    E  = "Robert Greene"
    F  = "Strategies Of War"
    G  = "Focus Area"
    H  = 471
    I  = "NA"
    J  = "Rs. 799"
    K  = 1
    L  = "1-1-D"
    M  = "Paperback"
    N  = 2007
    O  = "English"
    P  = "Self Help"
    Q  = "Strategy"
    R  = "No"
    S  = "NF"
    V  = 0
    X  = "Male"
    Y  = "978-1-8619-7978-0"
    Z  = 1
    AA = "First Floor"
    AB = "Practical"
    AC = "Yes"
    AD = "Yes"
    AE = 10.0
    AF = "New"
}

Add-BookRow -RowNum 185 -Values @{
    A  = 184
    B  = "Book"
    C  = "The 50th Law"
    D  = "Yes"
    E  = "Robert Greene"
    F  = "Biography Of 50cent"
    G  = "Focus Area"
    H  = 291
    I  = "NA"
    J  = "Rs. 699"
    K  = 1
    L  = "1-1-D"
    M  = "Paperback"
    N  = 2013
    O  = "English"
    P  = "Self Help"
    Q  = "Biography"
    R  = "No"
    S  = "NF"
    V  = 0
    X  = "Male"
    Y  = "978-1-8466-8079-3"
    Z  = 1
    AA = "First Floor"
    AB = "Practical"
    AC = "Yes"
    AD = "Yes"
    AE = 9.9
    AF = "New"
}

Add-BookRow -RowNum 186 -Values @{
    A  = 185
    B  = "Book"
    C  = "The Almanack Of Naval Ravikant"
    D  = "Yes"
    E  = "Eric Jorgenson"
    F  = "Insights From Angel Investor Naval Ravikanth"
    G  = "Focus Area"
    H  = 241
    I  = "Shane Parrish"
    J  = "Rs. 299"
    K  = 1
    L  = "1-1-D"
    M  = "Paperback"
    N  = 2020
    O  = "English"
    P  = "Self Help"
    Q  = "Strategy"
    R  = "No"
    S  = "NF"
    V  = 0
    X  = "Male"
    Y  = "978-93-5489-389-6"
    Z  = 1
    AA = "First Floor"
    AB = "Practical"
    AC = "Yes"
    AD = "Yes"
    AE = 10.0
    AF = "New"
}

"Added rows 183-186"
